$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the new "through" date
$wb.Worksheets.Item(1).Name = "Through 2021-11-06"

# Update the header label for the current-month column (B)
$ws.Range("B1").Value = "November 2021 (through November 06)"

# Row 3 - Garfield Park
$ws.Range("M3").Value = 4
$ws.Range("X3").Value = 1
$ws.Range("BE3").Value = 1

# Row 4 - Austin
$ws.Range("B4").Value = 2
$ws.Range("M4").Value = 3
$ws.Range("X4").Value = 1
$ws.Range("AT4").Value = 3
$ws.Range("BE4").Value = 1

# Row 5 - Humboldt Park
$ws.Range("B5").Value = 2
$ws.Range("M5").Value = 2

# Row 6 - West Town
$ws.Range("B6").Value = 2

# Row 7 - Englewood
$ws.Range("M7").Value = 2

# Row 11 - New City
$ws.Range("M11").Value = 2

# Row 17 - West Loop
$ws.Range("B17").Value = 2

# Row 20 - Lake View
$ws.Range("M20").Value = 1

# Row 21 - West Pullman
$ws.Range("AI21").Value = 1

# Row 28 - Uptown
$ws.Range("M28").Value = 1

# Row 36 - South Chicago
$ws.Range("B36").Value = 1

# Row 37 - Avalon Park
$ws.Range("BP37").Value = 1

# Row 39 - West Elsdon
$ws.Range("BP39").Value = 1

# Row 43 - Ashburn
$ws.Range("B43").Value = 2

# Row 58 - Fuller Park
$ws.Range("BE58").Value = 2

# Row 64 - Bridgeport
$ws.Range("B64").Value = 1

# Row 65 - Brighton Park
$ws.Range("BE65").Value = 2

# Row 66 - Chicago Lawn
$ws.Range("B66").Value = 4

# Row 68 - Douglas
$ws.Range("B68").Value = 1

# Row 84 - Morgan Park
$ws.Range("B84").Value = 3
